# Second editing pass: wording / phrasing touch-ups throughout the paper.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for:" $find
    }
    return $ok
}

# 1. "This raise follow-up questions:" -> "This raises the following questions:"
Replace-Text "raise follow-up questions" "raises the following questions"

# 2. "...Market exists where all knowledge is instantly available [GoBack] to all parties. Thus" -> "...Market exists with all knowledge instantly available to all parties. Thus"
Replace-Text "exists where all knowledge is instantly available to all parties" "exists with all knowledge instantly available to all parties"

# 3. "...random walk, this is not a limitation." -> "...random walk; this is not a limitation."
Replace-Text "truly is a random walk, this is not a limitation" "truly is a random walk; this is not a limitation"

# 4. "Rudnick describes several algorithms" -> "Rudnick provides several algorithms"
Replace-Text "Rudnick describes several algorithms" "Rudnick provides several algorithms"

# 5. "...modeling discounting the existing vast commoditization..." -> "...modeling overly discounts the vast existing commoditization..."
Replace-Text "modeling discounting the existing vast commoditization" "modeling overly discounts the vast existing commoditization"

# 6a. "applied to the financial markets" -> "applied to financial markets"
Replace-Text "applied to the financial markets" "applied to financial markets"
# 6b. "the price or volume represent" -> "the price and volume represent"
Replace-Text "the price or volume represent" "the price and volume represent"
# 6c. remove curly quotes around "moving average" -> "may lead to a more efficient moving average"
Replace-Text ([string][char]0x201C + "moving average" + [string][char]0x201D) "moving average"

# 7. "...last 100. One approach..." -> "...last 100 years. One approach..."
Replace-Text "for the last 100. One approach" "for the last 100 years. One approach"

# 8. "accomplished through correlations of multiple feeds" -> "accomplished through correlating multiple feeds"
Replace-Text "accomplished through correlations of multiple feeds" "accomplished through correlating multiple feeds"

# 9. "(2) always acts rationally" -> "(2) everyone always acts rationally"
Replace-Text "(2) always acts rational" "(2) everyone always acts rational"

# 10. "...could be created, and that it could potentially beat the market..." -> "...could exist that beats the market..."
Replace-Text "could be created, and that it could potentially beat the market" "could exist that beats the market"

# 11. "1. Preparing Next Quote Frame" -> "1. Preparing the Quote Frame"
Replace-Text "Preparing Next Quote Frame" "Preparing the Quote Frame"

# 12. "...moving averages to account for this." -> "...moving averages to account for this scenario."
Replace-Text "moving averages to account for this." "moving averages to account for this scenario."

# 13a. "...which are often simple aggregates." -> "...which are often simple aggregates such as net positive volume." + paragraph split
Replace-Text "which are often simple aggregates." "which are often simple aggregates such as net positive volume.`r"

# 13b. "...move in opposite directions, perhaps a sharp move..." -> "...move in opposite directions. Perhaps a sharp move..." + paragraph split before "3. Determine"
Replace-Text "move in opposite directions, perhaps a sharp move in bonds acts as a signal to change the desired equity position. " "move in opposite directions. Perhaps a sharp move in bonds acts as a signal to change the desired equity position.`r"

# 13c. "(alpha)" -> "(Alpha)"
Replace-Text "profitability (alpha) in financial markets" "profitability (Alpha) in financial markets"

# 14. "Beta is the deviation of correlated assets" -> "Beta occurs through the deviation of correlated assets"
Replace-Text "Beta is the deviation of correlated assets" "Beta occurs through the deviation of correlated assets"

# 15. "Vega is 1% increase" -> "Vega is gained by a 1% increase"
Replace-Text "Vega is 1% increase" "Vega is gained by a 1% increase"

# 16. "Gamma is increase rate of delta" -> "Gamma is the increase rate of Delta"
Replace-Text "Gamma is increase rate of delta" "Gamma is the increase rate of Delta"

# 17. "...delta gained per 1$ increase; and Theta is cost (interest) of 1 day passing." -> "...Delta gained per 1$ increase of the asset; and Theta is the cost (interest) from one day passing."
Replace-Text "gained per 1$ increase; and Theta is cost (interest) of 1 day passing" "gained per 1$ increase of the asset; and Theta is the cost (interest) from one day passing"

# 18. "(0, +Delta, -Vega, -Gamma, 0)" -> "(0, +Delta, -Vega, +Gamma, 0)"
Replace-Text "(0, +Delta, -Vega, -Gamma, 0)" "(0, +Delta, -Vega, +Gamma, 0)"

# 19. "If the current state exceeds the threshold of desired state;" -> "If the desired state exceeds the threshold relative to the current state;"
Replace-Text "If the current state exceeds the threshold of desired state;" "If the desired state exceeds the threshold relative to the current state;"

# 20. "cals, or (d) selling 1 x 30 Delta put" -> "cals, and (d) selling 1 x 30 Delta put"
Replace-Text "calls, or (d) selling 1 x 30 Delta put" "calls, and (d) selling 1 x 30 Delta put"

# 21. "The minimization of the other dimensions" -> "The minimization of other dimensions"
Replace-Text "The minimization of the other dimensions" "The minimization of other dimensions"
